$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "edad"
$ws.Range("D1").Value = "costo"

# Update data rows: column C becomes numeric "edad", column D becomes numeric "costo"
# (previously D held text "carrera" and E held numeric "costo")
$ws.Range("C2").Value = 17
$ws.Range("D2").Value = 5590000

$ws.Range("C3").Value = 18
$ws.Range("D3").Value = 5431212

$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 5202121

# Remove the now-unused column E entirely
$ws.Range("E1:E4").Delete()
